$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix J1 formula: relative R1 -> absolute $R$1 -----------------------
$ws.Range("J1").Formula = "=CONCATENATE(`$Q`$1,A1,`$S`$1,B1,`$S`$1,C1,`$S`$1,D1,`$S`$1,E1,`$S`$1,F1,`$R`$1)"

# --- 2. Add the four new user rows (19, 21, 23, 25) -------------------------
# Row 19 - Moira Parkel Vefill (B19 was typed with a leading quote-prefix
# apostrophe, so it keeps the quotePrefix cell style in the original file)
$ws.Cells.Item(19, 1).Value = 9043278
$ws.Cells.Item(19, 2).Value = "' 'Moira'"
$ws.Cells.Item(19, 3).Value = " 'Parkel'"
$ws.Cells.Item(19, 4).Value = " 'Vefill'"
$ws.Cells.Item(19, 5).Value = " 'parkelVef10@gmail.com'"
$ws.Cells.Item(19, 6).Value = 400
$ws.Cells.Item(19, 7).Value = 10

# Row 21 - Bobin Coronel Pesquero
$ws.Cells.Item(21, 1).Value = 4782107
$ws.Cells.Item(21, 2).Value = " 'Bobin'"
$ws.Cells.Item(21, 3).Value = " 'Coronel'"
$ws.Cells.Item(21, 4).Value = " 'Pesquero'"
$ws.Cells.Item(21, 5).Value = " 'nibobcorpes@hotmail.es'"
$ws.Cells.Item(21, 6).Value = 20
$ws.Cells.Item(21, 7).Value = 11

# Row 23 - Espinete Gonzalez Redondo
$ws.Cells.Item(23, 1).Value = 4309871
$ws.Cells.Item(23, 2).Value = " 'Espinete'"
$ws.Cells.Item(23, 3).Value = " 'Gonzalez'"
$ws.Cells.Item(23, 4).Value = " 'Redondo'"
$ws.Cells.Item(23, 5).Value = " 'redGonzes29@yahoo.es'"
$ws.Cells.Item(23, 6).Value = 3921
$ws.Cells.Item(23, 7).Value = 12

# Row 25 - Marta Guzman Fernandez (email text is missing its closing quote
# in the source data, reproduced verbatim)
$ws.Cells.Item(25, 1).Value = 3421987
$ws.Cells.Item(25, 2).Value = " 'Marta'"
$ws.Cells.Item(25, 3).Value = " 'Guzmán'"
$ws.Cells.Item(25, 4).Value = " 'Fernandez'"
$ws.Cells.Item(25, 5).Value = " 'martaGuzFer300297@gmail.com"
$ws.Cells.Item(25, 6).Value = 98
$ws.Cells.Item(25, 7).Value = 13

# --- 3. Rebuild the J-column "insert into" formulas for every data row -----
# The K-column reference used to be empty, truncating the generated SQL, so
# every data row formula is switched to use the absolute $R$1 terminator
# (which holds the ");" string) instead.
$dataRows = @(3, 5, 7, 9, 11, 13, 15, 17, 19, 21, 23, 25)
foreach ($r in $dataRows) {
    $formula = "=CONCATENATE(`$Q`$1,A$r,`$S`$1,B$r,`$S`$1,C$r,`$S`$1,D$r,`$S`$1,E$r,`$S`$1,F$r,`$R`$1)"
    $ws.Range("J${r}:J25").Formula = $formula
}

# Writing a formula across J<r>:J25 also stamps the intervening blank rows
# (even rows, which hold no other data) with shared-formula member cells.
# Those rows do not exist in the source file, so clear them back out again;
# an empty row with no cells left is dropped from the saved XML entirely.
$blankRows = @(4, 6, 8, 10, 12, 14, 16, 18, 20, 22, 24)
foreach ($r in $blankRows) {
    $ws.Cells.Item($r, 10).ClearContents()
}

# --- 4. Restore the active cell selection shown in the saved file ----------
$ws.Range("F26").Select()
